# Rename the AHB-Diff header row from the generic "_old" / "_new" suffixes
# to the concrete format-version suffixes "_FV2410" / "_FV2504", then turn
# the header+data range into a real Excel Table (ListObject) and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (A1:U1) ------------------------------------
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn A1:U81 into an Excel Table (ListObject) -----------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U81"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (row 1) ---------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
